# Updates cryptos list figures (price/volume columns + two coin-rank swaps)
# to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that sometimes looks numeric (e.g. "0.9997").
# Prefixing the value with a leading apostrophe forces Excel to keep it as
# text, exactly as typing it manually would - this preserves trailing
# zeros (e.g. "2.090") and avoids scientific notation (e.g. "0.000007671").
function Set-TextCell([string]$addr, [string]$text) {
    $ws.Range($addr).Value = "'" + $text
}

$ws.Range('D2').Value = '30.597.01'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '1.919.99'
$ws.Range('E3').Value = '  -0.17%  '
Set-TextCell 'D4' '0.9997'
$ws.Range('E4').Value = '  -0.03%  '
Set-TextCell 'D5' '245.28'
$ws.Range('E5').Value = '  -0.82%  '
Set-TextCell 'D6' '0.9999'
Set-TextCell 'D7' '0.4827'
$ws.Range('E7').Value = '  +1.72%  '
Set-TextCell 'D8' '0.2898'
$ws.Range('E8').Value = '  -0.35%  '
Set-TextCell 'D9' '0.06816'
$ws.Range('E9').Value = '  -0.24%  '
Set-TextCell 'D10' '111.98'
$ws.Range('E10').Value = '  +5.92%  '
Set-TextCell 'D11' '19.49'
$ws.Range('E11').Value = '  +5.80%  '
$ws.Range('D12').Value = '1.915.00'
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D13' '5.488'
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 'D14' '0.07572'
$ws.Range('E14').Value = '  -1.69%  '
Set-TextCell 'D15' '0.6734'
$ws.Range('E15').Value = '  +0.24%  '
Set-TextCell 'D16' '295.47'
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('D17').Value = '30.571.88'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 'D18' '0.000007671'
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 'D19' '13.03'
$ws.Range('E19').Value = '  +0.53%  '
Set-TextCell 'D20' '0.9998'
Set-TextCell 'D21' '5.521'
$ws.Range('E21').Value = '  -0.87%  '
$ws.Range('D22').Value = '2.163.53'
$ws.Range('E22').Value = '  -0.58%  '
Set-TextCell 'D23' '0.9995'
$ws.Range('E23').Value = '  -0.05%  '
Set-TextCell 'D24' '6.445'
$ws.Range('E24').Value = '  -0.60%  '
Set-TextCell 'D25' '9.494'
$ws.Range('E25').Value = '  -0.07%  '
Set-TextCell 'D26' '166.94'
$ws.Range('E26').Value = '  -0.39%  '
Set-TextCell 'D27' '20.35'
$ws.Range('E27').Value = '  -3.19%  '
Set-TextCell 'D28' '2.090'
$ws.Range('E28').Value = '  -1.69%  '
Set-TextCell 'D29' '0.1065'
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('E30').Value = '  +2.87%  '
Set-TextCell 'D31' '4.134'
$ws.Range('E31').Value = '  -1.13%  '
Set-TextCell 'D32' '4.060'
$ws.Range('E32').Value = '  +0.15%  '
Set-TextCell 'D33' '0.04982'
$ws.Range('E33').Value = '  -0.69%  '
Set-TextCell 'D34' '0.7346'
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('E35').Value = '  -0.56%  '
$ws.Range('E36').Value = '  -0.42%  '
$ws.Range('E37').Value = '  -2.14%  '
Set-TextCell 'D38' '2.684'
$ws.Range('E38').Value = '  +0.07%  '
Set-TextCell 'D39' '2.026'
$ws.Range('E39').Value = '  -0.79%  '
Set-TextCell 'D40' '109.39'
$ws.Range('E40').Value = '  -2.13%  '
Set-TextCell 'D41' '0.4436'
$ws.Range('E41').Value = '  +0.60%  '
Set-TextCell 'D42' '0.8697'
$ws.Range('E42').Value = '  -0.37%  '
Set-TextCell 'D43' '5.860'
$ws.Range('E43').Value = '  -0.84%  '
Set-TextCell 'D45' '69.47'
$ws.Range('E45').Value = '  +2.38%  '
Set-TextCell 'D46' '7.252'
$ws.Range('E46').Value = '  -0.59%  '
Set-TextCell 'D47' '48.92'
$ws.Range('E47').Value = '  +1.16%  '
Set-TextCell 'D48' '9.212'
$ws.Range('E48').Value = '  -1.47%  '
Set-TextCell 'D49' '0.1230'
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell 'D50' '34.85'
$ws.Range('E50').Value = '  -0.50%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextCell 'D51' '0.2508'
$ws.Range('E51').Value = '  +0.02%  '
